$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = 1
$ws.Range("C2").Value = "aditya"
$ws.Range("D2").Value = "employee"
$ws.Range("E2").Value = "manufacturing"
$ws.Range("B2").Value = "C1001"

$ws.Range("B2").Select()
